$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update text values (translated strings)
$ws.Range("A2").Value = "Vase"
$ws.Range("B2").Value = "Cat"
$ws.Range("C2").Value = "Ukraine"

# Update numeric values
$ws.Range("D2").Value = 6
$ws.Range("E2").Value = 120
$ws.Range("F2").Value = 240
$ws.Range("F3").Value = 240
